# Update countries & provincias Spain
# Applies the 1-Apr-2020 10:20 -> 10:50 COVID data refresh:
#  - Re-sorted a few countries whose "Casos totales" (col B) overtook a
#    neighbour, so the row now shows the country whose total is appropriate
#    for that rank; the refreshed per-row stats follow the country.
#  - A handful of other rows get refreshed stats without any re-sort.
#  - The "Datos actualizados a ..." timestamp text is bumped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner, cell A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Abril de 2020 a las 10:50"

# row -> (country, totalCasos, nuevosCasos, casosActivos, recuperados, casosCriticos, muertesHoy, muertes)
$rows = @(
    @{ R = 4;   Pais = "Estados Unidos";               B = 188592; C = 62; D = 7251; E = 177285; F = 4576; G = 3;  H = 4056 },
    @{ R = 34;  Pais = "Rumania";                       B = 2245;   C = 0;  D = 220;  E = 1940;   F = 62;   G = 3;  H = 85 },
    @{ R = 72;  Pais = "Bosnia y Herzegovina";           B = 448;    C = 28; D = 17;   E = 418;    F = 1;    G = 0;  H = 13 },
    @{ R = 73;  Pais = "Letonia";                        B = 446;    C = 48; D = 1;    E = 445;    F = 3;    G = 0;  H = 0 },
    @{ R = 105; Pais = "Mauricio";                       B = 147;    C = 4;  D = 0;    E = 142;    F = 1;    G = 0;  H = 5 },
    @{ R = 109; Pais = "Estado de Palestina";            B = 134;    C = 15; D = 18;   E = 115;    F = 0;    G = 0;  H = 1 },
    @{ R = 110; Pais = "Brunei";                         B = 129;    C = 0;  D = 45;   E = 83;     F = 3;    G = 0;  H = 1 },
    @{ R = 111; Pais = "Martinica";                      B = 128;    C = 0;  D = 27;   E = 98;     F = 15;   G = 0;  H = 3 },
    @{ R = 112; Pais = "Montenegro";                     B = 120;    C = 11; D = 0;    E = 118;    F = 4;    G = 0;  H = 2 },
    @{ R = 146; Pais = "Etiopia";                        B = 29;     C = 3;  D = 2;    E = 27;     F = 2;    G = 0;  H = 0 },
    @{ R = 156; Pais = "San Martin (Parte Francesa)";    B = 16;     C = 1;  D = 2;    E = 13;     F = 0;    G = 0;  H = 1 },
    @{ R = 157; Pais = "Eritrea";                        B = 15;     C = 0;  D = 0;    E = 15;     F = 0;    G = 0;  H = 0 },
    @{ R = 158; Pais = "Guinea Ecuatorial";               B = 15;     C = 0;  D = 1;    E = 14;     F = 0;    G = 0;  H = 0 },
    @{ R = 159; Pais = "Birmania";                       B = 15;     C = 0;  D = 0;    E = 14;     F = 0;    G = 0;  H = 1 },
    @{ R = 160; Pais = "Bahamas";                        B = 15;     C = 1;  D = 1;    E = 14;     F = 0;    G = 0;  H = 0 }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Cells.Item($r, 1).Value = $row.Pais
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
}
